# Apply the table style change described by the commit diff.
#
# Slide 16 contains a single table (graphicFrame, the 3rd shape on the
# slide) whose table style was switched from the deck's custom
# "Table_0" style ({BE7478D8-88B0-4556-9401-DE1BE05C4624}, defined in
# ppt/tableStyles.xml) to PowerPoint's built-in "No Style, Table Grid"
# style ({9EA1BBEB-E5AF-4EBC-98B6-76FC43C712BD}).
#
# Table styles can't be assigned through the Table.Style property
# directly (PowerPoint raises "Table styles cannot be assigned through
# a property - call Table.ApplyStyle(...) instead"), so ApplyStyle is
# used with the built-in style's GUID.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table
$tbl.ApplyStyle("{9EA1BBEB-E5AF-4EBC-98B6-76FC43C712BD}")
